$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "59.326.67"
Set-TextValue "E2" "  +0.57%  "
Set-TextValue "D3" "2.638.74"
Set-TextValue "E3" "  -0.51%  "
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "530.35"
Set-TextValue "E5" "  +1.55%  "
Set-TextValue "D6" "145.79"
Set-TextValue "E6" "  +0.74%  "
Set-TextValue "E7" "  -0.08%  "
Set-TextValue "D8" "0.571"
Set-TextValue "E8" "  -0.48%  "
Set-TextValue "E9" "  -3.69%  "
Set-TextValue "E10" "  +1.03%  "
Set-TextValue "E11" "  +0.42%  "
Set-TextValue "E12" "  +0.45%  "
Set-TextValue "D13" "3.107.79"
Set-TextValue "E13" "  +0.09%  "
Set-TextValue "D14" "59.710.52"
Set-TextValue "E14" "  +1.17%  "
Set-TextValue "D15" "20.77"
Set-TextValue "E15" "  -1.90%  "
Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000136"
Set-TextValue "E16" "  -0.02%  "
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "2.630.77"
Set-TextValue "E17" "  -1.71%  "
Set-TextValue "D18" "344.19"
Set-TextValue "E18" "  +1.07%  "
Set-TextValue "D19" "4.43"
Set-TextValue "E19" "  +1.03%  "
Set-TextValue "D20" "10.56"
Set-TextValue "E20" "  +2.02%  "
Set-TextValue "D21" "6.36"
Set-TextValue "E21" "  +0.24%  "
Set-TextValue "E22" "  +0.02%  "
Set-TextValue "D23" "66.22"
Set-TextValue "E23" "  +4.07%  "
Set-TextValue "D24" "0.415"
Set-TextValue "E24" "  +0.69%  "
Set-TextValue "E25" "  +1.66%  "
Set-TextValue "D26" "2.767.84"
Set-TextValue "E26" "  -0.14%  "
Set-TextValue "E27" "  -0.17%  "
Set-TextValue "D28" "7.17"
Set-TextValue "E28" "  +1.07%  "
Set-TextValue "D29" "0.0₃0796"
Set-TextValue "E29" "  -0.72%  "
Set-TextValue "E30" "  -0.06%  "
Set-TextValue "D31" "6.34"
Set-TextValue "E31" "  -4.39%  "
Set-TextValue "D32" "1.61"
Set-TextValue "E32" "  +1.43%  "
Set-TextValue "D33" "18.98"
Set-TextValue "E33" "  +1.11%  "
Set-TextValue "D34" "149.68"
Set-TextValue "E34" "  +0.47%  "
Set-TextValue "D35" "4.16"
Set-TextValue "E35" "  -0.17%  "
Set-TextValue "D36" "1.17"
Set-TextValue "E36" "  -1.85%  "
Set-TextValue "D37" "0.856"
Set-TextValue "E37" "  -4.42%  "
Set-TextValue "D38" "0.849"
Set-TextValue "E38" "  -3.69%  "
Set-TextValue "D39" "36.44"
Set-TextValue "E39" "  -0.44%  "
Set-TextValue "D40" "1.46"
Set-TextValue "E40" "  -1.54%  "
Set-TextValue "D41" "3.61"
Set-TextValue "E41" "  +0.58%  "
Set-TextValue "D42" "0.997"
Set-TextValue "E42" "  -0.02%  "
Set-TextValue "D43" "0.0979"
Set-TextValue "E44" "  -2.57%  "
Set-TextValue "D45" "268.83"
Set-TextValue "E45" "  -2.54%  "
Set-TextValue "B46" "WhiteBITCoin"
Set-TextValue "C46" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D46" "10.72"
Set-TextValue "E46" "  +1.85%  "
Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "19.26"
Set-TextValue "E47" "  -3.18%  "
Set-TextValue "D48" "0.0532"
Set-TextValue "E48" "  -0.61%  "
Set-TextValue "D49" "2.035.49"
Set-TextValue "E49" "  +0.15%  "
Set-TextValue "B50" "VeChain"
Set-TextValue "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D50" "0.0229"
Set-TextValue "E50" "  +0.44%  "
Set-TextValue "B51" "RenderToken"
Set-TextValue "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "4.68"
Set-TextValue "E51" "  -2.17%  "
